# "Generate Report for Archive"
#
# This reproduces a re-generation of the localization status report:
#   1. The handoff/localization "Status" for the d40532e7... item moves on
#      from "Ready for handoff" to "In Translation". That status string is
#      shared by the Overview sheet (columns E/F, row 2) and by the per
#      language detail sheets ("zh-cn"/"de-de", column C, row 2), so all of
#      those cells are updated together.
#   2. The two "date" columns on the Overview sheet (E:F) and the matching
#      "Status" column (C) on each per-language sheet are narrower in the
#      refreshed report, reflecting the regenerated column widths.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Update the status text everywhere it appears.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# 2. Narrow the regenerated columns (was ~17.22 chars, now ~13.41 chars).
$newWidth = 12.42

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth
